$wb = $excel.ActiveWorkbook

# The two sheets that contain this event data: "展览" (Exhibitions) and
# "全部类型" (All types). Both mirror rows 2-4 with the same F-column
# ("想去人数" / want-to-go count) values that need updating.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 4861
    $ws.Range("F3").Value = 144
    $ws.Range("F4").Value = 848
}
